# Apply cryptocurrency price / volume(1h) updates to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.055.11"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "1.793.46"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.51%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "228.17"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.42%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "31.21"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.18"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.281"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.67%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0662"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "2.050.19"
$ws.Range("E13").Value = "  -0.31%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "11.30"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +10.37%  "
$ws.Range("D15").Value = "1.793.96"
$ws.Range("E15").Value = "  -0.02%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.636"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "34.093.10"
$ws.Range("E17").Value = "  -1.15%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.22"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.70%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "69.72"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "253.77"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.99%  "
$ws.Range("D21").Value = "0.0₃0745"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("E22").Value = "  +0.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.29"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("E25").Value = "  -1.41%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.98"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.36%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.63"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.115"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.87%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  +0.40%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.91"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0518"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "1.488.36"
$ws.Range("E36").Value = "  -5.54%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.06"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.635"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.00%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0188"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "83.88"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.51%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  -0.25%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.907"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  -3.91%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").Value = "1.947.85"
$ws.Range("E47").Value = "  -0.06%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.71"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("E49").Value = "  +0.37%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "11.82"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "51.48"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.90%  "
